# Adds the "deployment phase" dialog strings to both the "en" and "es"
# sheets of the language workbook, then leaves the "es" sheet active
# (matching the commit's final cursor/tab state).

$wb = $excel.ActiveWorkbook
$wsEn = $wb.Worksheets.Item("en")
$wsEs = $wb.Worksheets.Item("es")

# --- "en" sheet: new rows 71-76 (key in column A, English text in column B) ---
$enRows = @(
    @("deploy_quad_0",  "Before we begin, we must verify Robert's coordinates."),
    @("deploy_quad_1",  "First, you must choose the correct quadrant."),
    @("deploy_coord_0", "Excellent! Now we need to input the actual coordinate numbers."),
    @("deploy_coord_1", "Use the numpad to input the X and Y coordinates."),
    @("deploy_coord_2", "You can switch which coordinates to input by pressing the left or right arrow."),
    @("deploy_coord_3", "Once you are satisfied with the coordinate numbers, press the VERIFY button to proceed.")
)

$row = 71
foreach ($pair in $enRows) {
    $wsEn.Cells.Item($row, 1).Value = $pair[0]
    $cellB = $wsEn.Cells.Item($row, 2)
    $cellB.Value = $pair[1]
    $cellB.VerticalAlignment = -4108  # xlCenter - matches other wrapped/translated cells in column B
    $row = $row + 1
}

# --- "es" sheet: new rows 68-76, keys only (translations not done yet) ---
$esKeys = @(
    "deployment",
    "verify",
    "out_of_bounds",
    "deploy_quad_0",
    "deploy_quad_1",
    "deploy_coord_0",
    "deploy_coord_1",
    "deploy_coord_2",
    "deploy_coord_3"
)

$row = 68
foreach ($key in $esKeys) {
    $wsEs.Cells.Item($row, 1).Value = $key
    $row = $row + 1
}

# --- View state: selections on each sheet ---
$wsEn.Range("A68:A76").Select()
$wsEs.Range("B66").Select()

# "es" ends up as the active/selected tab
$wsEs.Activate()
